$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values are treated as text so numeric-looking
# strings (e.g. "211.22") are not auto-converted to numbers, then restore
# the default cell style so no visual/style change is introduced.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.650.44"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "1.599.42"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "211.22"
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("D6").Value = "0.518"
$ws.Range("E6").Value = "  +1.02%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -1.36%  "
$ws.Range("D10").Value = "19.42"
$ws.Range("E10").Value = "  -1.51%  "
$ws.Range("D11").Value = "0.0837"
$ws.Range("E11").Value = "  +0.21%  "
$ws.Range("D12").Value = "1.820.27"
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").Value = "1.601.20"
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").Value = "0.522"
$ws.Range("E15").Value = "  -1.00%  "
$ws.Range("D16").Value = "64.79"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("D17").Value = "26.615.93"
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("D18").Value = "0.0₃0731"
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "208.46"
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "1.00"
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("D21").Value = "6.93"
$ws.Range("E21").Value = "  +2.34%  "
$ws.Range("D22").Value = "4.26"
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("E23").Value = "  -3.15%  "
$ws.Range("D24").Value = "8.86"
$ws.Range("D25").Value = "145.57"
$ws.Range("E25").Value = "  -0.99%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "7.17"
$ws.Range("E27").Value = "  -1.32%  "
$ws.Range("E28").Value = "  +0.39%  "
$ws.Range("D29").Value = "15.32"
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  -0.41%  "
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("D32").Value = "3.24"
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("E33").Value = "  -0.78%  "
$ws.Range("D34").Value = "2.92"
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("D35").Value = "1.282.74"
$ws.Range("E35").Value = "  -2.61%  "
$ws.Range("E36").Value = "  +1.68%  "
$ws.Range("E37").Value = "  -0.30%  "
$ws.Range("E38").Value = "  -0.39%  "
$ws.Range("D39").Value = "0.842"
$ws.Range("E39").Value = "  +1.53%  "
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").Value = "5.46"
$ws.Range("E41").Value = "  +0.90%  "
$ws.Range("E42").Value = "  +1.12%  "
$ws.Range("D43").Value = "0.786"
$ws.Range("E43").Value = "  -0.49%  "
$ws.Range("D44").Value = "63.86"
$ws.Range("E44").Value = "  +0.70%  "
$ws.Range("E45").Value = "  +9.60%  "
$ws.Range("D46").Value = "1.733.26"
$ws.Range("E46").Value = "  +0.29%  "
$ws.Range("D47").Value = "89.78"
$ws.Range("E47").Value = "  -0.34%  "
$ws.Range("E48").Value = "  -0.86%  "
$ws.Range("D49").Value = "0.0₆0105"
$ws.Range("E49").Value = "  -1.11%  "
$ws.Range("E50").Value = "  +3.56%  "
$ws.Range("D51").Value = "0.0506"
$ws.Range("E51").Value = "  -0.98%  "

$ws.Range("D2:D51").Style = "Normal"

